$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Fix the placeholder project numbers in column Q (rows 4-9) which all
# incorrectly shared the same value "0000-1000" - give them unique values.
$ws.Range("Q4").Value = "0000-1001"
$ws.Range("Q5").Value = "0000-1002"
$ws.Range("Q6").Value = "0000-1003"
$ws.Range("Q7").Value = "0000-1004"
$ws.Range("Q8").Value = "0000-1005"
$ws.Range("Q9").Value = "0000-1006"

# Update the frozen-pane view state: scroll the bottom-right pane so that
# column J becomes the top-left visible column, and select Q10.
$ws.Activate()
$excel.ActiveWindow.TopLeftCell = "J3"
$ws.Range("Q10").Select()
